$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "244.26"
Set-TextValue $ws.Range("E2") "-0.52%"
Set-TextValue $ws.Range("G2") "10"
Set-TextValue $ws.Range("D3") "27.16"
Set-TextValue $ws.Range("E3") "3.01%"
Set-TextValue $ws.Range("G3") "10"
Set-TextValue $ws.Range("D4") "5.150"
Set-TextValue $ws.Range("E4") "0.20%"
Set-TextValue $ws.Range("G4") "10"
Set-TextValue $ws.Range("E5") "0.61%"
Set-TextValue $ws.Range("G5") "10"
Set-TextValue $ws.Range("D6") "6.474"
Set-TextValue $ws.Range("E6") "-0.25%"
Set-TextValue $ws.Range("G6") "10"
Set-TextValue $ws.Range("D7") "0.8161"
Set-TextValue $ws.Range("E7") "0.01%"
Set-TextValue $ws.Range("G7") "10"
Set-TextValue $ws.Range("D8") "0.8310"
Set-TextValue $ws.Range("E8") "-1.94%"
Set-TextValue $ws.Range("G8") "10"
Set-TextValue $ws.Range("D9") "0.1326"
Set-TextValue $ws.Range("E9") "-0.99%"
Set-TextValue $ws.Range("G9") "10"
Set-TextValue $ws.Range("D10") "0.06898"
Set-TextValue $ws.Range("E10") "-0.70%"
Set-TextValue $ws.Range("G10") "10"
Set-TextValue $ws.Range("D11") "0.02888"
Set-TextValue $ws.Range("E11") "1.22%"
Set-TextValue $ws.Range("G11") "10"
Set-TextValue $ws.Range("D12") "0.09389"
Set-TextValue $ws.Range("E12") "-0.16%"
Set-TextValue $ws.Range("G12") "10"
Set-TextValue $ws.Range("E13") "-1.10%"
Set-TextValue $ws.Range("G13") "10"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D14") "0.006185"
Set-TextValue $ws.Range("E14") "-0.98%"
Set-TextValue $ws.Range("G14") "10"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D15") "3.606"
Set-TextValue $ws.Range("E15") "1.62%"
Set-TextValue $ws.Range("G15") "10"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D16") "3.020"
Set-TextValue $ws.Range("E16") "-0.01%"
Set-TextValue $ws.Range("G16") "10"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D17") "2.307"
Set-TextValue $ws.Range("E17") "8.91%"
Set-TextValue $ws.Range("G17") "10"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.009987"
Set-TextValue $ws.Range("E18") "1,561.34%"
Set-TextValue $ws.Range("G18") "10"
Set-TextValue $ws.Range("D19") "0.3113"
Set-TextValue $ws.Range("E19") "-0.70%"
Set-TextValue $ws.Range("G19") "10"
Set-TextValue $ws.Range("D20") "0.03113"
Set-TextValue $ws.Range("E20") "-3.31%"
Set-TextValue $ws.Range("G20") "10"
Set-TextValue $ws.Range("E21") "-2.18%"
Set-TextValue $ws.Range("G21") "10"
Set-TextValue $ws.Range("D22") "3.732"
Set-TextValue $ws.Range("E22") "-0.67%"
Set-TextValue $ws.Range("G22") "10"
Set-TextValue $ws.Range("D23") "0.04485"
Set-TextValue $ws.Range("E23") "-4.58%"
Set-TextValue $ws.Range("G23") "10"
Set-TextValue $ws.Range("G24") "10"
Set-TextValue $ws.Range("E25") "-1.94%"
Set-TextValue $ws.Range("G25") "10"
Set-TextValue $ws.Range("D26") "0.004491"
Set-TextValue $ws.Range("E26") "-2.48%"
Set-TextValue $ws.Range("G26") "10"
Set-TextValue $ws.Range("D27") "0.00009791"
Set-TextValue $ws.Range("E27") "1.97%"
Set-TextValue $ws.Range("G27") "10"
Set-TextValue $ws.Range("D28") "0.0001394"
Set-TextValue $ws.Range("E28") "0.31%"
Set-TextValue $ws.Range("G28") "10"
Set-TextValue $ws.Range("G29") "10"
Set-TextValue $ws.Range("G30") "10"
Set-TextValue $ws.Range("G31") "10"
Set-TextValue $ws.Range("G32") "10"
Set-TextValue $ws.Range("G33") "10"
Set-TextValue $ws.Range("G34") "10"
Set-TextValue $ws.Range("G35") "10"
Set-TextValue $ws.Range("G36") "10"
Set-TextValue $ws.Range("G37") "10"
Set-TextValue $ws.Range("G38") "10"
Set-TextValue $ws.Range("G39") "10"
Set-TextValue $ws.Range("D40") "0.03642"
Set-TextValue $ws.Range("E40") "-0.31%"
Set-TextValue $ws.Range("G40") "10"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006045"
Set-TextValue $ws.Range("E41") "-1.42%"
Set-TextValue $ws.Range("G41") "10"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1051"
Set-TextValue $ws.Range("E42") "-0.25%"
Set-TextValue $ws.Range("G42") "10"
Set-TextValue $ws.Range("D43") "0.002571"
Set-TextValue $ws.Range("E43") "0.48%"
Set-TextValue $ws.Range("G43") "10"
Set-TextValue $ws.Range("D44") "0.008190"
Set-TextValue $ws.Range("E44") "3.35%"
Set-TextValue $ws.Range("G44") "10"
Set-TextValue $ws.Range("D45") "0.00005314"
Set-TextValue $ws.Range("E45") "0.07%"
Set-TextValue $ws.Range("G45") "10"
Set-TextValue $ws.Range("D46") "0.00000000749"
Set-TextValue $ws.Range("E46") "-0.11%"
Set-TextValue $ws.Range("G46") "10"
Set-TextValue $ws.Range("E47") "-18.45%"
Set-TextValue $ws.Range("G47") "10"
Set-TextValue $ws.Range("D48") "0.002601"
Set-TextValue $ws.Range("E48") "26.96%"
Set-TextValue $ws.Range("G48") "10"
Set-TextValue $ws.Range("D49") "0.00002098"
Set-TextValue $ws.Range("E49") "-0.11%"
Set-TextValue $ws.Range("G49") "10"
Set-TextValue $ws.Range("D50") "0.0001998"
Set-TextValue $ws.Range("E50") "-0.11%"
Set-TextValue $ws.Range("G50") "10"
Set-TextValue $ws.Range("G51") "10"
